$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 118
$ws.Range("H118").Value = 4246.943
$ws.Range("I118").Value = 483
$ws.Range("J118").Value = 9892.857
$ws.Range("K118").Value = 1449
$ws.Range("L118").Value = 29678.571
$ws.Range("M118").Value = 208
$ws.Range("N118").Value = -32992.571

# Row 138
$ws.Range("H138").Value = 2834.4285
$ws.Range("I138").Value = 1291.9354
$ws.Range("K138").Value = 3875.8062
$ws.Range("M138").Value = 1264.1938

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 244.14285
$ws.Range("I5").Value = 279.66666
$ws.Range("J5").Value = 217.5
$ws.Range("K5").Value = 279.66666
$ws.Range("L5").Value = 217.5
$ws.Range("M5").Value = -167.66666
$ws.Range("N5").Value = -441.5

# Row 9
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

# Row 10
$ws.Range("H10").Value = 10002.5
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 10002.5
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 10002.5
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -10342.5

# Row 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

# Row 20
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

# Row 22
$ws.Range("H22").Value = 4079.7334
$ws.Range("I22").Value = 626
$ws.Range("J22").Value = 10987.2
$ws.Range("K22").Value = 626
$ws.Range("L22").Value = 10987.2
$ws.Range("M22").Value = -327
$ws.Range("N22").Value = -11585.2

# Row 32
$ws.Range("H32").Value = 7501.11
$ws.Range("I32").Value = 5319.325
$ws.Range("J32").Value = 16228.25
$ws.Range("K32").Value = 5319.325
$ws.Range("L32").Value = 16228.25
$ws.Range("M32").Value = -5032.325
$ws.Range("N32").Value = -16802.25

# Row 61
$ws.Range("H61").Value = 2527.1936
$ws.Range("I61").Value = 2309
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2309
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -2097
$ws.Range("N61").Value = -4424

# Row 74
$ws.Range("H74").Value = 234755.44
$ws.Range("I74").Value = 2025.1578
$ws.Range("J74").Value = 2003505.6
$ws.Range("K74").Value = 2025.1578
$ws.Range("L74").Value = 2003505.6
$ws.Range("M74").Value = -1151.1578
$ws.Range("N74").Value = -2005253.6

# Row 77
$ws.Range("H77").Value = 234755.44
$ws.Range("I77").Value = 2025.1578
$ws.Range("J77").Value = 2003505.6
$ws.Range("K77").Value = 10125.789
$ws.Range("L77").Value = 10017528
$ws.Range("M77").Value = -5757.789000000001
$ws.Range("N77").Value = -10026264

# Row 132
$ws.Range("H132").Value = 24145.152
$ws.Range("I132").Value = 32967.938
$ws.Range("J132").Value = 3978.7856
$ws.Range("K132").Value = 98903.81400000001
$ws.Range("L132").Value = 11936.3568
$ws.Range("M132").Value = -96373.81400000001
$ws.Range("N132").Value = -16996.3568

# Row 133
$ws.Range("H133").Value = 34482.715
$ws.Range("J133").Value = 34482.715
$ws.Range("L133").Value = 34482.715
$ws.Range("N133").Value = -39542.715

# Row 136
$ws.Range("H136").Value = 2527.1936
$ws.Range("I136").Value = 2309
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 6927
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -4377
$ws.Range("N136").Value = -17100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 244.14285
$ws.Range("I4").Value = 279.66666
$ws.Range("J4").Value = 217.5
$ws.Range("K4").Value = 279.66666
$ws.Range("L4").Value = 217.5
$ws.Range("M4").Value = -164.66666
$ws.Range("N4").Value = -447.5

# Row 14
$ws.Range("H14").Value = 1580
$ws.Range("I14").Value = 1580
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1580
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1408
$ws.Range("N14").ClearContents()

# Row 15
$ws.Range("H15").Value = 3500
$ws.Range("J15").Value = 3500
$ws.Range("L15").Value = 3500
$ws.Range("N15").Value = -3954

# Row 16
$ws.Range("H16").Value = 10000
$ws.Range("J16").Value = 10000
$ws.Range("L16").Value = 10000
$ws.Range("N16").Value = -10340

# Row 22
$ws.Range("H22").Value = 816.125
$ws.Range("I22").Value = 589.8
$ws.Range("J22").Value = 1193.3334
$ws.Range("K22").Value = 589.8
$ws.Range("L22").Value = 1193.3334
$ws.Range("M22").Value = -416.8
$ws.Range("N22").Value = -1539.3334

# Row 80
$ws.Range("H80").Value = 148.75
$ws.Range("J80").Value = 164.92857
$ws.Range("L80").Value = 164.92857
$ws.Range("N80").Value = -2160.92857

# Row 83
$ws.Range("H83").Value = 148.75
$ws.Range("J83").Value = 164.92857
$ws.Range("L83").Value = 824.6428500000001
$ws.Range("N83").Value = -10808.64285

# Row 107
$ws.Range("H107").Value = 4301.5483
$ws.Range("I107").Value = 5755
$ws.Range("J107").Value = 1249.3
$ws.Range("K107").Value = 5755
$ws.Range("L107").Value = 1249.3
$ws.Range("M107").Value = -3835
$ws.Range("N107").Value = -5089.3

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2515.054
$ws.Range("I31").Value = 1772.4
$ws.Range("J31").Value = 2790.111
$ws.Range("K31").Value = 1772.4
$ws.Range("L31").Value = 2790.111
$ws.Range("M31").Value = -1477.4
$ws.Range("N31").Value = -3380.111

# Row 34
$ws.Range("H34").Value = 2515.054
$ws.Range("I34").Value = 1772.4
$ws.Range("J34").Value = 2790.111
$ws.Range("K34").Value = 1772.4
$ws.Range("L34").Value = 2790.111
$ws.Range("M34").Value = -1570.4
$ws.Range("N34").Value = -3194.111

# Row 36
$ws.Range("H36").Value = 3000
$ws.Range("I36").Value = 3000
$ws.Range("K36").Value = 3000
$ws.Range("M36").Value = -2612

# Row 40
$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 3000
$ws.Range("K40").Value = 3000
$ws.Range("M40").Value = -2840

# Row 132
$ws.Range("H132").Value = 1630.8
$ws.Range("I132").Value = 1362
$ws.Range("J132").Value = 3378
$ws.Range("K132").Value = 4086
$ws.Range("L132").Value = 10134
$ws.Range("M132").Value = -1556
$ws.Range("N132").Value = -15194

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 165012.58
$ws.Range("I2").Value = 396012
$ws.Range("J2").Value = 13
$ws.Range("K2").Value = 2376072
$ws.Range("L2").Value = 78
$ws.Range("M2").Value = -2375959
$ws.Range("N2").Value = -304

# Row 10
$ws.Range("H10").Value = 410
$ws.Range("I10").Value = 30
$ws.Range("J10").Value = 600
$ws.Range("K10").Value = 90
$ws.Range("L10").Value = 1800
$ws.Range("M10").Value = 49
$ws.Range("N10").Value = -2078

# Row 26
$ws.Range("H26").Value = 40230.1
$ws.Range("I26").Value = 425.25
$ws.Range("J26").Value = 66766.664
$ws.Range("K26").Value = 1275.75
$ws.Range("L26").Value = 200299.992
$ws.Range("M26").Value = -987.75
$ws.Range("N26").Value = -200875.992

# Row 32
$ws.Range("H32").Value = 1687.2195
$ws.Range("I32").Value = 172.63637
$ws.Range("J32").Value = 2242.5667
$ws.Range("K32").Value = 517.9091100000001
$ws.Range("L32").Value = 6727.7001
$ws.Range("M32").Value = -234.9091100000001
$ws.Range("N32").Value = -7293.7001

# Row 33
$ws.Range("H33").Value = 6864.6665
$ws.Range("I33").Value = 325
$ws.Range("J33").Value = 7870.769
$ws.Range("K33").Value = 1950
$ws.Range("L33").Value = 47224.614
$ws.Range("M33").Value = -1667
$ws.Range("N33").Value = -47790.614

# Row 38
$ws.Range("H38").Value = 100.958336
$ws.Range("I38").Value = 114.94118
$ws.Range("J38").Value = 67
$ws.Range("K38").Value = 344.82354
$ws.Range("L38").Value = 201
$ws.Range("M38").Value = 2.17646000000002
$ws.Range("N38").Value = -895

# Row 39
$ws.Range("H39").Value = 2746.647
$ws.Range("J39").Value = 2946
$ws.Range("L39").Value = 8838
$ws.Range("N39").Value = -9426

# Row 44
$ws.Range("H44").Value = 5166.6665
$ws.Range("I44").Value = 3000
$ws.Range("J44").Value = 6714.2856
$ws.Range("K44").Value = 9000
$ws.Range("L44").Value = 20142.8568
$ws.Range("M44").Value = -8602
$ws.Range("N44").Value = -20938.8568

# Row 46
$ws.Range("H46").Value = 2382.6
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 2502.7856
$ws.Range("K46").Value = 2100
$ws.Range("L46").Value = 7508.3568
$ws.Range("M46").Value = -2009
$ws.Range("N46").Value = -7690.3568

# Row 50
$ws.Range("H50").Value = 108.14286
$ws.Range("I50").Value = 31.6
$ws.Range("J50").Value = 299.5
$ws.Range("K50").Value = 94.80000000000001
$ws.Range("L50").Value = 898.5
$ws.Range("M50").Value = 386.2
$ws.Range("N50").Value = -1860.5

# Row 51
$ws.Range("H51").Value = 3799
$ws.Range("I51").Value = 850
$ws.Range("J51").Value = 4192.2
$ws.Range("K51").Value = 2550
$ws.Range("L51").Value = 12576.6
$ws.Range("M51").Value = -2090
$ws.Range("N51").Value = -13496.6

# Row 53
$ws.Range("H53").Value = 108.14286
$ws.Range("I53").Value = 31.6
$ws.Range("J53").Value = 299.5
$ws.Range("K53").Value = 94.80000000000001
$ws.Range("L53").Value = 898.5
$ws.Range("M53").Value = 386.2
$ws.Range("N53").Value = -1860.5

# Row 57
$ws.Range("H57").Value = 1400
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

# Row 58
$ws.Range("H58").Value = 4278.778
$ws.Range("I58").Value = 3000
$ws.Range("J58").Value = 4354
$ws.Range("K58").Value = 9000
$ws.Range("L58").Value = 13062
$ws.Range("M58").Value = -8872
$ws.Range("N58").Value = -13318

# Row 121
$ws.Range("H121").Value = 1173.6207
$ws.Range("I121").Value = 341.15384
$ws.Range("J121").Value = 1850
$ws.Range("K121").Value = 1023.46152
$ws.Range("L121").Value = 5550
$ws.Range("M121").Value = 286.5384799999999
$ws.Range("N121").Value = -8170

# Row 131
$ws.Range("H131").Value = 6667410
$ws.Range("J131").Value = 1243.25
$ws.Range("L131").Value = 3729.75
$ws.Range("N131").Value = -13809.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2352.5918
$ws.Range("I132").Value = 1871.5758
$ws.Range("J132").Value = 3344.6875
$ws.Range("K132").Value = 5614.7274
$ws.Range("L132").Value = 10034.0625
$ws.Range("M132").Value = -3084.7274
$ws.Range("N132").Value = -15094.0625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 2999.1738
$ws.Range("I136").Value = 633.2195
$ws.Range("J136").Value = 22400
$ws.Range("K136").Value = 1899.6585
$ws.Range("L136").Value = 67200
$ws.Range("M136").Value = 650.3415
$ws.Range("N136").Value = -72300
